# Apply updated "dSF" (column F) values for specific rows on Sheet1.
# Row number in the sheet = data-index + 2 (row 1 is header).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 2
    4  = 0
    5  = -1
    8  = 4
    9  = -3
    10 = -2
    12 = -10
    14 = -1
    15 = 1
    21 = -2
    28 = 2
    30 = -1
    33 = 0
    34 = 1
    35 = -1
    36 = 0
    40 = -3
    44 = 9
    47 = -1
    49 = -3
    50 = 4
    71 = 3
    74 = 4
    75 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
